# Upgrade description labels in 外观检查.xlsx
# Rename header columns:
#   损坏类型 -> 缺损类型
#   病害描述 -> 缺损描述
# across every worksheet (all three sheets share the same header layout).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $d1 = $ws.Range("D1")
    if ($d1.Text -eq "损坏类型") {
        $d1.Value = "缺损类型"
    }
    $e1 = $ws.Range("E1")
    if ($e1.Text -eq "病害描述") {
        $e1.Value = "缺损描述"
    }
    # Reset the active selection back to A1 so the saved sheetView has no
    # explicit <selection> override (matches Excel's default-state omission).
    $ws.Range("A1").Select()
}

$wb.Worksheets.Item(1).Activate()
$wb.Worksheets.Item(1).Range("A1").Select()
